# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# cells for the 6b361d66... entry on both the zh-cn and de-de sheets, as a result of a
# freshly generated handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 01:02:43"
$wsZhCn.Range("H2").Value = "2016-03-25 01:03:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 01:02:47"
$wsDeDe.Range("H2").Value = "2016-03-25 01:03:19"
